# Updates the "cryptos" price list with freshly scraped Price/Volume(1h)
# figures (and, for two swapped-rank pairs, the Coin/Link too).
#
# The "Price" column (D) holds numbers formatted like "27.011.74" or
# "0.07230" that must stay exactly as literal text (multiple '.' separators,
# leading/trailing zeros, etc. are meaningful and must survive unchanged).
# Assigning such a string straight to .Value would make Excel auto-detect it
# as a number and mangle it (e.g. drop the leading/trailing zeros or render
# it in scientific notation), so those values are written with a leading
# apostrophe, exactly like a user typing '0.07229 into a cell, which forces
# Excel to keep it as quoted text instead of re-interpreting it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.004.62"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").Value = "'1.849.49"
$ws.Range("E3").Value = "  +1.25%  "
$ws.Range("D4").Value = "'1.014"
$ws.Range("E4").Value = "  +0.57%  "
$ws.Range("E5").Value = "  +0.44%  "
$ws.Range("D6").Value = "'309.75"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").Value = "'0.4772"
$ws.Range("E7").Value = "  +1.97%  "
$ws.Range("D8").Value = "'0.3679"
$ws.Range("E8").Value = "  +2.10%  "
$ws.Range("D9").Value = "'0.07229"
$ws.Range("E9").Value = "  +1.38%  "
$ws.Range("E10").Value = "  +3.37%  "
$ws.Range("D11").Value = "'19.73"
$ws.Range("E11").Value = "  +1.65%  "
$ws.Range("D12").Value = "'0.07721"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.859.82"
$ws.Range("E13").Value = "  +1.35%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'5.329"
$ws.Range("E14").Value = "  +1.08%  "
$ws.Range("D15").Value = "'6.435"
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("D16").Value = "'88.93"
$ws.Range("E16").Value = "  +1.56%  "
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").Value = "'0.000008638"
$ws.Range("E18").Value = "  +1.08%  "
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("D20").Value = "'27.045.88"
$ws.Range("D21").Value = "'14.50"
$ws.Range("E21").Value = "  +1.96%  "
$ws.Range("D22").Value = "'5.056"
$ws.Range("E22").Value = "  +0.67%  "
$ws.Range("E23").Value = "  +0.76%  "
$ws.Range("D24").Value = "'1.928"
$ws.Range("E24").Value = "  +1.43%  "
$ws.Range("D25").Value = "'152.65"
$ws.Range("E25").Value = "  -0.24%  "
$ws.Range("D26").Value = "'18.22"
$ws.Range("E26").Value = "  +1.65%  "
$ws.Range("D27").Value = "'2.004"
$ws.Range("E27").Value = "  +1.23%  "
$ws.Range("D28").Value = "'114.49"
$ws.Range("E28").Value = "  +0.60%  "
$ws.Range("D29").Value = "'5.002"
$ws.Range("E29").Value = "  +2.75%  "
$ws.Range("D30").Value = "'0.08901"
$ws.Range("E30").Value = "  +1.04%  "
$ws.Range("D31").Value = "'3.318"
$ws.Range("E31").Value = "  +5.58%  "
$ws.Range("E32").Value = "  +1.18%  "
$ws.Range("D33").Value = "'0.7463"
$ws.Range("E33").Value = "  +1.27%  "
$ws.Range("D34").Value = "'4.503"
$ws.Range("E34").Value = "  +1.44%  "
$ws.Range("D35").Value = "'2.724"
$ws.Range("E35").Value = "  -3.21%  "
$ws.Range("D36").Value = "'1.113"
$ws.Range("E36").Value = "  +3.18%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.05276"
$ws.Range("E37").Value = "  +2.27%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.01954"
$ws.Range("E38").Value = "  +1.09%  "
$ws.Range("D39").Value = "'2.977"
$ws.Range("E39").Value = "  +2.38%  "
$ws.Range("D40").Value = "'0.5241"
$ws.Range("E40").Value = "  +3.75%  "
$ws.Range("D41").Value = "'7.012"
$ws.Range("E41").Value = "  +1.91%  "
$ws.Range("E42").Value = "  +1.14%  "
$ws.Range("D43").Value = "'8.217"
$ws.Range("E43").Value = "  +2.03%  "
$ws.Range("E44").Value = "  +5.98%  "
$ws.Range("E45").Value = "  +1.86%  "
$ws.Range("E46").Value = "  +0.44%  "
$ws.Range("D47").Value = "'101.56"
$ws.Range("E47").Value = "  +3.82%  "
$ws.Range("D48").Value = "'1.611"
$ws.Range("E48").Value = "  +2.47%  "
$ws.Range("D49").Value = "'65.47"
$ws.Range("E49").Value = "  +2.30%  "
$ws.Range("D50").Value = "'0.06060"
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("D51").Value = "'0.8879"
$ws.Range("E51").Value = "  +4.14%  "
